$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in row 2 data values ---
# Columns: A=ID, B=ServerID, C=Name, D=MaxOnline, E=CpuCount, F=IP, G=Port
$ws.Range("A2").Value = "WorldServer_1"

$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "000103001"

$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "WorldServer_1"

$ws.Range("D2").Value = 5000
$ws.Range("E2").Value = 1

$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "127.0.0.1"

$ws.Range("G2").Value = 3001

# --- Column width changes ---
$ws.Columns.Item(2).ColumnWidth = 13.714285714285714
$ws.Columns.Item(3).ColumnWidth = 17.857142857142858
$ws.Columns.Item(5).ColumnWidth = 11.285714285714286
$ws.Columns.Item(6).ColumnWidth = 15.0
$ws.Columns.Item(7).ColumnWidth = 10.285714285714286

# --- Selection moves to G1 ---
$ws.Range("G1").Select()

# --- Remove data validations from the sheet ---
$ws.Cells.Validation.Delete()
